# "расписание на вторник.xlsx" — correct the misspelled teacher name
# "Юрьевнга" -> "Юрьевна" (Моргуненко Елена Юрьевна) in the row-4 lesson
# cells, and move the active selection to H7 (matching the saved view
# state in the authored workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the "проект" cell first, then the "английский" cells, so the
# corrected shared strings are (re)created in the same order as in the
# target workbook (проект, then английский, appended at the end of the
# shared-string table).
$ws.Range("H4").Value = "10Т проект Моргуненко Елена Юрьевна"
$ws.Range("F4").Value = "10Т английский Моргуненко Елена Юрьевна"
$ws.Range("G4").Value = "10Т английский Моргуненко Елена Юрьевна"

# Update the saved selection/active cell to H7.
$ws.Range("H7").Select()
